$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.559.50'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '2.609.11'
$ws.Range("E3").Value = '  +0.90%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.61'
$ws.Range("E5").Value = '  +2.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.77'
$ws.Range("E6").Value = '  +1.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.567'
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.102'
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.335'
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").Value = '3.068.87'
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").Value = '59.484.09'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.68'
$ws.Range("E15").Value = '  +1.33%  '
$ws.Range("D16").Value = '2.619.20'
$ws.Range("E16").Value = '  +1.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000133'
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '341.60'
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.36'
$ws.Range("E19").Value = '  +1.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.11'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.39'
$ws.Range("E21").Value = '  -2.00%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.44'
$ws.Range("E23").Value = '  +2.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.409'
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.165'
$ws.Range("E25").Value = '  -1.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.24'
$ws.Range("E27").Value = '  +3.26%  '
$ws.Range("D28").Value = '0.0₃0749'
$ws.Range("E28").Value = '  +3.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.67'
$ws.Range("E30").Value = '  +5.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.82'
$ws.Range("E31").Value = '  -2.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.85'
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.64'
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.98'
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("E35").Value = '  -0.36%  '
$ws.Range("B36").Value = 'Stacks'
$ws.Range("C36").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.47'
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("B37").Value = 'SuiNetwork'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.840'
$ws.Range("E37").Value = '  +3.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.832'
$ws.Range("E38").Value = '  +0.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.54'
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '274.89'
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.599'
$ws.Range("E42").Value = '  +1.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.72'
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0955'
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0524'
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("D46").Value = '1.950.18'
$ws.Range("E46").Value = '  -1.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.53'
$ws.Range("E47").Value = '  +3.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0223'
$ws.Range("E48").Value = '  +1.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.52'
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.15'
$ws.Range("E50").Value = '  -2.22%  '
$ws.Range("E51").Value = '  +0.52%  '

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
